$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1100
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 1400
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 1400
$ws.Range("M40").Value = -825
$ws.Range("N40").Value = -1750

$ws.Range("H64").Value = 3000
$ws.Range("J64").Value = 3333.3333
$ws.Range("L64").Value = 3333.3333
$ws.Range("N64").Value = -3829.3333

$ws.Range("H67").Value = 3000
$ws.Range("J67").Value = 3333.3333
$ws.Range("L67").Value = 3333.3333
$ws.Range("N67").Value = -5049.3333

$ws.Range("H116").Value = 673426.8
$ws.Range("I116").Value = 1002639.6
$ws.Range("K116").Value = 1002639.6
$ws.Range("M116").Value = -999197.6

$ws.Range("H132").Value = 24489568
$ws.Range("I132").Value = 25101432
$ws.Range("K132").Value = 75304296
$ws.Range("M132").Value = -75301766

$ws.Range("H133").Value = 59642
$ws.Range("J133").Value = 59642
$ws.Range("L133").Value = 59642
$ws.Range("N133").Value = -69762

$ws.Range("H137").Value = 3804.8462
$ws.Range("I137").Value = 3792.4375
$ws.Range("K137").Value = 11377.3125
$ws.Range("M137").Value = -8827.3125

$ws.Range("H138").Value = 4127.8667
$ws.Range("I138").Value = 2309.7778
$ws.Range("J138").Value = 4375.788
$ws.Range("K138").Value = 6929.3334
$ws.Range("L138").Value = 13127.364
$ws.Range("M138").Value = -1789.3334
$ws.Range("N138").Value = -23407.364

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1053.0869
$ws.Range("I2").Value = 1040.0555
$ws.Range("J2").Value = 1100
$ws.Range("K2").Value = 1040.0555
$ws.Range("L2").Value = 1100
$ws.Range("M2").Value = -927.0554999999999
$ws.Range("N2").Value = -1326

$ws.Range("H80").Value = 34555.453
$ws.Range("J80").Value = 34555.453
$ws.Range("L80").Value = 34555.453
$ws.Range("N80").Value = -36551.453

$ws.Range("H83").Value = 34555.453
$ws.Range("J83").Value = 34555.453
$ws.Range("L83").Value = 103666.359
$ws.Range("N83").Value = -113650.359

$ws.Range("H111").Value = 34644
$ws.Range("J111").Value = 34644
$ws.Range("L111").Value = 34644
$ws.Range("N111").Value = -42824

$ws.Range("H116").Value = 1053.0869
$ws.Range("I116").Value = 1040.0555
$ws.Range("J116").Value = 1100
$ws.Range("K116").Value = 1040.0555
$ws.Range("L116").Value = 1100
$ws.Range("M116").Value = 1253.9445
$ws.Range("N116").Value = -5688

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1053.0869
$ws.Range("I3").Value = 1040.0555
$ws.Range("J3").Value = 1100
$ws.Range("K3").Value = 1040.0555
$ws.Range("L3").Value = 1100
$ws.Range("M3").Value = -926.0554999999999
$ws.Range("N3").Value = -1328

$ws.Range("H129").Value = 45334.75
$ws.Range("J129").Value = 45334.75
$ws.Range("L129").Value = 45334.75
$ws.Range("N129").Value = -55334.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5105.472
$ws.Range("I31").Value = 3140
$ws.Range("J31").Value = 5422.484
$ws.Range("K31").Value = 3140
$ws.Range("L31").Value = 5422.484
$ws.Range("M31").Value = -2845
$ws.Range("N31").Value = -6012.484

$ws.Range("H34").Value = 5105.472
$ws.Range("I34").Value = 3140
$ws.Range("J34").Value = 5422.484
$ws.Range("K34").Value = 3140
$ws.Range("L34").Value = 5422.484
$ws.Range("M34").Value = -2938
$ws.Range("N34").Value = -5826.484

$ws.Range("H68").Value = 48439.383
$ws.Range("J68").Value = 48439.383
$ws.Range("L68").Value = 48439.383
$ws.Range("N68").Value = -49937.383

$ws.Range("H71").Value = 48439.383
$ws.Range("J71").Value = 48439.383
$ws.Range("L71").Value = 145318.149
$ws.Range("N71").Value = -152806.149

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H110").Value = 40780
$ws.Range("J110").Value = 40780
$ws.Range("L110").Value = 40780
$ws.Range("N110").Value = -48960

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 5998.5
$ws.Range("J95").Value = 5998.5
$ws.Range("L95").Value = 17995.5
$ws.Range("N95").Value = -22113.5

$ws.Range("H113").Value = 605.2973
$ws.Range("I113").Value = 638.44446
$ws.Range("J113").Value = 573.8946999999999
$ws.Range("K113").Value = 1915.33338
$ws.Range("L113").Value = 1721.6841
$ws.Range("M113").Value = 254.66662
$ws.Range("N113").Value = -6061.6841

$ws.Range("H114").Value = 1542
$ws.Range("I114").Value = 455.75
$ws.Range("J114").Value = 2266.1667
$ws.Range("K114").Value = 1367.25
$ws.Range("L114").Value = 6798.500100000001
$ws.Range("M114").Value = 1886.75
$ws.Range("N114").Value = -13306.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 20061
$ws.Range("I43").Value = 2500
$ws.Range("K43").Value = 2500
$ws.Range("M43").Value = -2349

$ws.Range("H57").Value = 32372.5
$ws.Range("I57").Value = 29500
$ws.Range("J57").Value = 33330
$ws.Range("K57").Value = 29500
$ws.Range("L57").Value = 33330
$ws.Range("M57").Value = -28680
$ws.Range("N57").Value = -34970

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H135").Value = 65455.645
$ws.Range("I135").Value = 139500
$ws.Range("J135").Value = 59759.92
$ws.Range("K135").Value = 139500
$ws.Range("L135").Value = 59759.92
$ws.Range("M135").Value = -134430
$ws.Range("N135").Value = -69899.92

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 17858138
$ws.Range("I81").Value = 20090294
$ws.Range("J81").Value = 900
$ws.Range("K81").Value = 40180588
$ws.Range("L81").Value = 1800
$ws.Range("M81").Value = -40179527
$ws.Range("N81").Value = -3922

$ws.Range("H84").Value = 17858138
$ws.Range("I84").Value = 20090294
$ws.Range("J84").Value = 900
$ws.Range("K84").Value = 200902940
$ws.Range("L84").Value = 9000
$ws.Range("M84").Value = -200897636
$ws.Range("N84").Value = -19608

$ws.Range("H107").Value = 700.7368
$ws.Range("J107").Value = 1520
$ws.Range("L107").Value = 4560
$ws.Range("N107").Value = -8400

$ws.Range("H132").Value = 5558931
$ws.Range("I132").Value = 3885.8572
$ws.Range("J132").Value = 13335994
$ws.Range("K132").Value = 11657.5716
$ws.Range("L132").Value = 40007982
$ws.Range("M132").Value = -9127.571599999999
$ws.Range("N132").Value = -40013042

$ws.Range("H136").Value = 5212.4707
$ws.Range("I136").Value = 1788.875
$ws.Range("J136").Value = 8255.666999999999
$ws.Range("K136").Value = 5366.625
$ws.Range("L136").Value = 24767.001
$ws.Range("M136").Value = -2816.625
$ws.Range("N136").Value = -29867.001
